$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row data: row number -> new Price (column D) and new Volume(1h) (column E).
# $null for D means the Price column is unchanged for that row.
$updates = @(
    [PSCustomObject]@{ Row = 2; D = '89.713.09'; E = '  -0.80%  ' }
    [PSCustomObject]@{ Row = 3; D = '3.070.08'; E = '  -1.78%  ' }
    [PSCustomObject]@{ Row = 4; D = '1.00'; E = '  -0.30%  ' }
    [PSCustomObject]@{ Row = 5; D = $null; E = '  +9.81%  ' }
    [PSCustomObject]@{ Row = 6; D = '616.76'; E = '  -0.64%  ' }
    [PSCustomObject]@{ Row = 7; D = $null; E = '  -5.53%  ' }
    [PSCustomObject]@{ Row = 8; D = '0.366'; E = '  +1.24%  ' }
    [PSCustomObject]@{ Row = 9; D = '0.999'; E = '  -0.11%  ' }
    [PSCustomObject]@{ Row = 10; D = '3.067.12'; E = '  -1.86%  ' }
    [PSCustomObject]@{ Row = 11; D = '0.702'; E = '  -4.53%  ' }
    [PSCustomObject]@{ Row = 12; D = $null; E = '  -0.03%  ' }
    [PSCustomObject]@{ Row = 13; D = '0.0000248'; E = '  +1.24%  ' }
    [PSCustomObject]@{ Row = 14; D = '34.64'; E = '  -1.30%  ' }
    [PSCustomObject]@{ Row = 15; D = '89.369.99'; E = '  -0.97%  ' }
    [PSCustomObject]@{ Row = 16; D = '5.37'; E = '  -4.75%  ' }
    [PSCustomObject]@{ Row = 17; D = '3.637.47'; E = '  -2.01%  ' }
    [PSCustomObject]@{ Row = 18; D = '3.061.13'; E = '  -2.49%  ' }
    [PSCustomObject]@{ Row = 19; D = '3.80'; E = '  +1.53%  ' }
    [PSCustomObject]@{ Row = 20; D = '0.0000212'; E = '  +0.38%  ' }
    [PSCustomObject]@{ Row = 21; D = '13.74'; E = '  -5.07%  ' }
    [PSCustomObject]@{ Row = 22; D = '430.85'; E = '  -6.60%  ' }
    [PSCustomObject]@{ Row = 23; D = '5.40'; E = '  +0.81%  ' }
    [PSCustomObject]@{ Row = 24; D = '8.70'; E = '  -3.70%  ' }
    [PSCustomObject]@{ Row = 25; D = '5.55'; E = '  -3.57%  ' }
    [PSCustomObject]@{ Row = 26; D = '11.70'; E = '  -4.35%  ' }
    [PSCustomObject]@{ Row = 27; D = '81.53'; E = '  -13.97%  ' }
    [PSCustomObject]@{ Row = 28; D = $null; E = '  -2.23%  ' }
    [PSCustomObject]@{ Row = 30; D = $null; E = '  +41.37%  ' }
    [PSCustomObject]@{ Row = 31; D = '0.159'; E = '  -2.47%  ' }
    [PSCustomObject]@{ Row = 32; D = '8.95'; E = '  -2.71%  ' }
    [PSCustomObject]@{ Row = 33; D = $null; E = '  -9.66%  ' }
    [PSCustomObject]@{ Row = 34; D = '4.24'; E = '  +65.58%  ' }
    [PSCustomObject]@{ Row = 35; D = '25.57'; E = '  -3.71%  ' }
    [PSCustomObject]@{ Row = 36; D = '0.150'; E = '  +3.00%  ' }
    [PSCustomObject]@{ Row = 37; D = '7.11'; E = '  +1.83%  ' }
    [PSCustomObject]@{ Row = 38; D = '490.15'; E = '  -4.98%  ' }
    [PSCustomObject]@{ Row = 39; D = '3.61'; E = '  +0.79%  ' }
    [PSCustomObject]@{ Row = 40; D = $null; E = '  -2.96%  ' }
    [PSCustomObject]@{ Row = 41; D = $null; E = '  -5.06%  ' }
    [PSCustomObject]@{ Row = 42; D = '0.0897'; E = '  -2.61%  ' }
    [PSCustomObject]@{ Row = 43; D = $null; E = '  -0.68%  ' }
    [PSCustomObject]@{ Row = 45; D = '0.399'; E = '  -6.31%  ' }
    [PSCustomObject]@{ Row = 46; D = '156.43'; E = '  +3.88%  ' }
    [PSCustomObject]@{ Row = 47; D = $null; E = '  -6.78%  ' }
    [PSCustomObject]@{ Row = 48; D = '0.671'; E = '  -6.80%  ' }
    [PSCustomObject]@{ Row = 49; D = '44.42'; E = '  -1.91%  ' }
    [PSCustomObject]@{ Row = 50; D = $null; E = '  -0.25%  ' }
    [PSCustomObject]@{ Row = 51; D = $null; E = '  -5.17%  ' }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        $dCell = $ws.Cells.Item($u.Row, 4)
        # Force text format so numeric-looking price strings (e.g. "1.00",
        # "89.713.09") are preserved exactly instead of being parsed as numbers.
        $dCell.NumberFormat = "@"
        $dCell.Value = $u.D
    }
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}
